$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("A2").Value = "A 11"
$ws.Range("B2").Value = "B 12"
$ws.Range("C2").Value = "C 13"
$ws.Range("D2").Value = "D 14"

$ws.Range("A3").Value = "A 21"
$ws.Range("B3").Value = "B 22"
$ws.Range("C3").Value = "C 23"
$ws.Range("D3").Value = "D 24"

$ws.Range("A4").Value = "A 31"
$ws.Range("B4").Value = "B 32"
$ws.Range("C4").Value = "C 33"
$ws.Range("D4").Value = "D 34"
